# Apply numeric updates per the commit diff, grouped by worksheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 51882.332
$ws.Range("J3").Value = 51882.332
$ws.Range("L3").Value = 51882.332
$ws.Range("N3").Value = -52110.332

$ws.Range("H40").Value = 2415.0312
$ws.Range("I40").Value = 2384.158
$ws.Range("J40").Value = 2460.1538
$ws.Range("K40").Value = 2384.158
$ws.Range("L40").Value = 2460.1538
$ws.Range("M40").Value = -2209.158
$ws.Range("N40").Value = -2810.1538

$ws.Range("H64").Value = 4273.8887
$ws.Range("I64").Value = 3633.3333
$ws.Range("K64").Value = 3633.3333
$ws.Range("M64").Value = -3385.3333

$ws.Range("H67").Value = 4273.8887
$ws.Range("I67").Value = 3633.3333
$ws.Range("K67").Value = 3633.3333
$ws.Range("M67").Value = -2775.3333

$ws.Range("H69").Value = 9007.5
$ws.Range("J69").Value = 10015
$ws.Range("L69").Value = 30045
$ws.Range("N69").Value = -31793

$ws.Range("H72").Value = 9007.5
$ws.Range("J72").Value = 10015
$ws.Range("L72").Value = 90135
$ws.Range("N72").Value = -98871

$ws.Range("H92").Value = 712.5172
$ws.Range("I92").Value = 592.6
$ws.Range("K92").Value = 592.6
$ws.Range("M92").Value = 655.4

$ws.Range("H100").Value = 43278.082
$ws.Range("I100").Value = 46723.863
$ws.Range("K100").Value = 46723.863
$ws.Range("M100").Value = -46182.863

$ws.Range("H102").Value = 51882.332
$ws.Range("J102").Value = 51882.332
$ws.Range("L102").Value = 51882.332
$ws.Range("N102").Value = -58372.332

$ws.Range("H107").Value = 75.5
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H115").Value = 406
$ws.Range("I115").Value = 406
$ws.Range("K115").Value = 1218
$ws.Range("M115").Value = 349

$ws.Range("H116").Value = 7338.1
$ws.Range("I116").Value = 4648
$ws.Range("J116").Value = 11373.25
$ws.Range("K116").Value = 4648
$ws.Range("L116").Value = 11373.25
$ws.Range("M116").Value = -1206
$ws.Range("N116").Value = -18257.25

$ws.Range("H127").Value = 3339.1155
$ws.Range("I127").Value = 1628.3334
$ws.Range("K127").Value = 4885.0002
$ws.Range("M127").Value = 74.9997999999996

$ws.Range("H137").Value = 11105.755
$ws.Range("I137").Value = 5422.7036
$ws.Range("K137").Value = 16268.1108
$ws.Range("M137").Value = -13718.1108

$ws.Range("H138").Value = 5369.5674
$ws.Range("J138").Value = 5285.1787
$ws.Range("L138").Value = 15855.5361
$ws.Range("N138").Value = -26135.5361

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10093.608
$ws.Range("I32").Value = 9223.282999999999
$ws.Range("J32").Value = 39249.5
$ws.Range("K32").Value = 9223.282999999999
$ws.Range("L32").Value = 39249.5
$ws.Range("M32").Value = -8936.282999999999
$ws.Range("N32").Value = -39823.5

$ws.Range("H74").Value = 8413.35
$ws.Range("I74").Value = 7345.6274
$ws.Range("K74").Value = 7345.6274
$ws.Range("M74").Value = -6471.6274

$ws.Range("H77").Value = 8413.35
$ws.Range("I77").Value = 7345.6274
$ws.Range("K77").Value = 36728.137
$ws.Range("M77").Value = -32360.137

$ws.Range("H122").Value = 2392.2827
$ws.Range("J122").Value = 3047.0908
$ws.Range("L122").Value = 9141.2724
$ws.Range("N122").Value = -14041.2724

$ws.Range("H132").Value = 2226.984
$ws.Range("I132").Value = 2141.0754
$ws.Range("K132").Value = 6423.226200000001
$ws.Range("M132").Value = -3893.226200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 286.66666
$ws.Range("J64").Value = 286.1
$ws.Range("L64").Value = 286.1
$ws.Range("N64").Value = -736.1

$ws.Range("H67").Value = 286.66666
$ws.Range("J67").Value = 286.1
$ws.Range("L67").Value = 286.1
$ws.Range("N67").Value = -1846.1

$ws.Range("H86").Value = 403413.1
$ws.Range("I86").Value = 835034.75
$ws.Range("J86").Value = 4993.077
$ws.Range("K86").Value = 835034.75
$ws.Range("L86").Value = 4993.077
$ws.Range("M86").Value = -833911.75
$ws.Range("N86").Value = -7239.077

$ws.Range("H89").Value = 403413.1
$ws.Range("I89").Value = 835034.75
$ws.Range("J89").Value = 4993.077
$ws.Range("K89").Value = 4175173.75
$ws.Range("L89").Value = 24965.385
$ws.Range("M89").Value = -4169557.75
$ws.Range("N89").Value = -36197.385

$ws.Range("H94").Value = 1228.44
$ws.Range("J94").Value = 261.75
$ws.Range("L94").Value = 261.75
$ws.Range("N94").Value = -1163.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4598.4653
$ws.Range("I58").Value = 3480.5312
$ws.Range("J58").Value = 5974.385
$ws.Range("K58").Value = 3480.5312
$ws.Range("L58").Value = 5974.385
$ws.Range("M58").Value = -3277.5312
$ws.Range("N58").Value = -6380.385

$ws.Range("H136").Value = 4598.4653
$ws.Range("I136").Value = 3480.5312
$ws.Range("J136").Value = 5974.385
$ws.Range("K136").Value = 10441.5936
$ws.Range("L136").Value = 17923.155
$ws.Range("M136").Value = -7891.5936
$ws.Range("N136").Value = -23023.155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2229.0264
$ws.Range("I5").Value = 1276.5264
$ws.Range("J5").Value = 3181.5264
$ws.Range("K5").Value = 3829.5792
$ws.Range("L5").Value = 9544.5792
$ws.Range("M5").Value = -3717.5792
$ws.Range("N5").Value = -9768.5792

$ws.Range("H76").Value = 4000.3333
$ws.Range("I76").Value = 4000.3333
$ws.Range("K76").Value = 12000.9999
$ws.Range("M76").Value = -11617.9999

$ws.Range("H79").Value = 4000.3333
$ws.Range("I79").Value = 4000.3333
$ws.Range("K79").Value = 12000.9999
$ws.Range("M79").Value = -10674.9999

$ws.Range("H113").Value = 2375.5625
$ws.Range("I113").Value = 2227.25
$ws.Range("K113").Value = 6681.75
$ws.Range("M113").Value = -4511.75

$ws.Range("H131").Value = 6674.9375
$ws.Range("I131").Value = 1621
$ws.Range("J131").Value = 7610.852
$ws.Range("K131").Value = 4863
$ws.Range("L131").Value = 22832.556
$ws.Range("M131").Value = 177
$ws.Range("N131").Value = -32912.556

$ws.Range("H135").Value = 2229.0264
$ws.Range("I135").Value = 1276.5264
$ws.Range("J135").Value = 3181.5264
$ws.Range("K135").Value = 11488.7376
$ws.Range("L135").Value = 28633.7376
$ws.Range("M135").Value = -8953.7376
$ws.Range("N135").Value = -33703.7376

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 44991.332
$ws.Range("I62").Value = 44991.332
$ws.Range("K62").Value = 44991.332
$ws.Range("M62").Value = -44305.332

$ws.Range("H65").Value = 44991.332
$ws.Range("I65").Value = 44991.332
$ws.Range("K65").Value = 134973.996
$ws.Range("M65").Value = -131541.996

$ws.Range("H107").Value = 444.23077
$ws.Range("I107").Value = 475.6
$ws.Range("K107").Value = 475.6
$ws.Range("M107").Value = 1444.4

$ws.Range("H132").Value = 3821.7322
$ws.Range("I132").Value = 4023.5881
$ws.Range("J132").Value = 1762.8
$ws.Range("K132").Value = 12070.7643
$ws.Range("L132").Value = 5288.4
$ws.Range("M132").Value = -9540.764299999999
$ws.Range("N132").Value = -10348.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 988.1818
$ws.Range("I22").Value = 697.1429000000001
$ws.Range("J22").Value = 1497.5
$ws.Range("K22").Value = 697.1429000000001
$ws.Range("L22").Value = 1497.5
$ws.Range("M22").Value = -402.1429000000001
$ws.Range("N22").Value = -2087.5

$ws.Range("H27").Value = 988.1818
$ws.Range("I27").Value = 697.1429000000001
$ws.Range("J27").Value = 1497.5
$ws.Range("K27").Value = 697.1429000000001
$ws.Range("L27").Value = 1497.5
$ws.Range("M27").Value = -590.1429000000001
$ws.Range("N27").Value = -1711.5

$ws.Range("H46").Value = 1695.0588
$ws.Range("I46").Value = 1008.9
$ws.Range("J46").Value = 1980.9584
$ws.Range("K46").Value = 1008.9
$ws.Range("L46").Value = 1980.9584
$ws.Range("M46").Value = -820.9
$ws.Range("N46").Value = -2356.9584

$ws.Range("H61").Value = 2189
$ws.Range("J61").Value = 2759.4
$ws.Range("L61").Value = 2759.4
$ws.Range("N61").Value = -3163.4

$ws.Range("H113").Value = 2189
$ws.Range("J113").Value = 2759.4
$ws.Range("L113").Value = 2759.4
$ws.Range("N113").Value = -7099.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 14174.5
$ws.Range("I58").Value = 14174.5
$ws.Range("K58").Value = 14174.5
$ws.Range("M58").Value = -13866.5

$ws.Range("H107").Value = 1706.1818
$ws.Range("I107").Value = 876.8
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 2630.4
$ws.Range("L107").Value = 30000
$ws.Range("M107").Value = -710.3999999999996
$ws.Range("N107").Value = -33840

$ws.Range("H113").Value = 681.75
$ws.Range("I113").Value = 651.6667
$ws.Range("J113").Value = 772
$ws.Range("K113").Value = 1955.0001
$ws.Range("L113").Value = 2316
$ws.Range("M113").Value = 214.9999
$ws.Range("N113").Value = -6656

$ws.Range("H122").Value = 5835
$ws.Range("I122").Value = 4703
$ws.Range("K122").Value = 14109
$ws.Range("M122").Value = -11659

$ws.Range("H132").Value = 124386.375
$ws.Range("I132").Value = 158529.67
$ws.Range("J132").Value = 31516.6
$ws.Range("K132").Value = 475589.01
$ws.Range("L132").Value = 94549.79999999999
$ws.Range("M132").Value = -473059.01
$ws.Range("N132").Value = -99609.79999999999

$ws.Range("H135").Value = 41998.832
$ws.Range("J135").Value = 41998.832
$ws.Range("L135").Value = 41998.832
$ws.Range("N135").Value = -52138.832

$ws.Range("H136").Value = 6820411.5
$ws.Range("I136").Value = 9678814
$ws.Range("J136").Value = 4221.3076
$ws.Range("K136").Value = 29036442
$ws.Range("L136").Value = 12663.9228
$ws.Range("M136").Value = -29033892
$ws.Range("N136").Value = -17763.9228
